$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet had a 2-row header (row1 = unit labels scattered across
# E/G/I/J/K, row2 = Hiver/Annee/Ete sub-labels across F..K) followed by the
# data table starting at row3. The new layout merges everything into a
# single header row (row1: idx/idx2/Name/Date Start/Date End/(m3/s)/(MW1)/
# (MW2)/(GWh) Winter/(GWh) Summer/(GWh) Year) with the data table shifting
# up to start at row2. Deleting the old row1 achieves exactly that shift
# while leaving the (already correctly styled) data rows untouched.
$ws.Rows("1").Delete()

# Build the new style used by the unit-header cells (F1:K1): same font as
# the existing "text" style (Arial 9) but keeping numFmtId General and no
# applyNumberFormat flag. Registering it as a transient named cell style and
# deleting the name afterwards leaves only the resulting cell-format record
# behind (matching how the workbook was produced), without leaving an extra
# named style / cellStyleXfs entry in the saved file.
$hdrStyleName = "__tmp_unit_header_style"
$hdrStyle = $wb.Styles.Add($hdrStyleName, "Normal")
$hdrStyle.Font.Name = "Arial"
$hdrStyle.Font.Size = 9

# Row 1: index / identification columns (A:E) - plain default style.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

# Row 1: unit columns (F:K) - use the new header style.
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

$ws.Range("F1:K1").Style = $hdrStyleName

# Drop the temporary named style now that its format has been stamped onto
# the cells; only the resulting cell format (cellXfs entry) should remain.
$wb.Styles.Item($hdrStyleName).Delete()

# Match the author's final selection (first data row highlighted).
$ws.Range("A2:K2").Select()
